$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-25 05:40:23"
$wsZhCn.Range("G4").Value = "2016-02-25 05:41:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-25 05:40:35"
$wsDeDe.Range("G4").Value = "2016-02-25 05:41:34"
